$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new problem entry as row 23
$ws.Range("A23").Value = "Largest Sum Contiguous Subarray (Kadane’s Algorithm)"
$ws.Range("B23").Value = "KadaneAlgorithm"

# Update the selected cell to match the new last row
$ws.Range("B23").Select()
